$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three rows were removed entirely from the species list:
#   old row 18 - "Scyliorhinus canicula" (1-RAP gear)
#   old row 27 - "Lophius budegassa"     (2-RAP gear)
#   old row 32 - "Scyliorhinus canicula" (2-RAP gear)
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(32).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(18).Delete()

# After the deletions the surviving rows shift up. Starting at the new
# row 32 (old row 35, "Anadara transversa") through the new last row 66
# (old row 69, "Wood NA") every record now carries an RF value in column I.
$ws.Range("I32:I66").Value = 82.75702127659575

# A handful of those same rows previously stored Numb = 0; those become -1.
$ws.Range("H37").Value = -1
$ws.Range("H38").Value = -1
$ws.Range("H40").Value = -1
$ws.Range("H60").Value = -1
$ws.Range("H63").Value = -1
$ws.Range("H66").Value = -1
